# Generate Report for Handoff
# Updates the localization-status report after a new handoff run:
#  - Priority for the 086775c0 file moves from "low" to "ht" (zh-cn & de-de)
#  - Latest Handoff Datetime for that file is refreshed (zh-cn)
#  - Latest HO Xliff Generate Date on the Overview sheet is refreshed

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# zh-cn sheet: Priority (col E) and Latest Handoff Datetime (col H) for rows 4-7
for ($r = 4; $r -le 7; $r++) {
    $zhcn.Cells.Item($r, 5).Value = "ht"
    $zhcn.Cells.Item($r, 8).Value = "2016-09-05 18:34:55"
}

# de-de sheet: Priority (col E) for rows 4-7.
# Latest Handoff Datetime (col H) shares the same underlying string as the
# Overview sheet's "Latest HO Xliff Generate Date" (both were the generic
# "2016-09-05 18:34:35" placeholder), so it also picks up the refreshed
# timestamp below.
for ($r = 4; $r -le 7; $r++) {
    $dede.Cells.Item($r, 5).Value = "ht"
    $dede.Cells.Item($r, 8).Value = "2016-09-05 18:35:01"
}

# Overview sheet: Latest HO Xliff Generate Date (col G) for rows 4-7
for ($r = 4; $r -le 7; $r++) {
    $overview.Cells.Item($r, 7).Value = "2016-09-05 18:35:01"
}
